# Fri, May 29, 2020  2:05:01 AM
#
# The deck's Design theme is switched from the custom "Integral" (Red
# Violet colour scheme) back to the default "Office Theme" colour
# scheme. Re-theming a deck like this also resets any table that was
# still wearing the old theme's custom table style back to PowerPoint's
# built-in "No Style, Table Grid" style.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Every table in the deck that still uses the old custom table style
#    ({DA455424-1A74-4819-831C-AEB893656D73}) picks up the built-in
#    "No Style, Table Grid" style instead.
# ---------------------------------------------------------------------
$oldStyle = "{DA455424-1A74-4819-831C-AEB893656D73}"
$newStyle = "{5792C567-EAAC-44BA-BE7A-3C9AA7888890}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldStyle) {
                $tbl.ApplyStyle($newStyle)
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Re-colour the presentation's theme from "Red Violet" to the
#    standard Office colour scheme.
# ---------------------------------------------------------------------
function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$master = $p.Designs.Item(1).SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

$officeColors = @(
    (RGBVal 0x00 0x00 0x00),  # dk1
    (RGBVal 0xFF 0xFF 0xFF),  # lt1
    (RGBVal 0x44 0x54 0x6A),  # dk2
    (RGBVal 0xE7 0xE6 0xE6),  # lt2
    (RGBVal 0x5B 0x9B 0xD5),  # accent1
    (RGBVal 0xED 0x7D 0x31),  # accent2
    (RGBVal 0xA5 0xA5 0xA5),  # accent3
    (RGBVal 0xFF 0xC0 0x00),  # accent4
    (RGBVal 0x44 0x72 0xC4),  # accent5
    (RGBVal 0x70 0xAD 0x47),  # accent6
    (RGBVal 0x05 0x63 0xC1),  # hlink
    (RGBVal 0x95 0x4F 0x72)   # folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i - 1]
}
